$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C header
$ws.Range("C1").Value = "pincode"

# Column C width
$ws.Columns("C").ColumnWidth = 16

# B4: empty cell that inherits the "data" style (font + vertical-center) plus a
# quote-prefix flag - achieved by copying B2's format onto B4, then toggling
# the quote-prefix (by entering and clearing a leading-apostrophe value).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Value = "'x"
$ws.Range("B4").ClearContents() | Out-Null

# C3: pincode value, entered with a leading apostrophe so Excel stores it as
# quote-prefixed text.
$ws.Range("C3").Value = "'5800  '"

# Update the selected cell to C3, matching the saved view state.
$ws.Range("C3").Select() | Out-Null
